# Updated cryptos list on Tue Oct  8 04:51:29 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) updates for each row.
# Each entry: row, newPrice text (or $null to leave unchanged), newVolume text (or $null to leave unchanged)
# NeedsTextFormat = $true for Price values that look like plain numbers, so Excel
# doesn't auto-convert them from text to a numeric value (the source data keeps
# these as text, e.g. "1.00", "26.90", preserving trailing zeros / exact formatting).
$updates = @(
    @{ Row = 2;  D = "62.620.44"; DText = $false; E = $null },
    @{ Row = 3;  D = "2.441.69";  DText = $false; E = "  -1.44%  " },
    @{ Row = 4;  D = "0.999";     DText = $true;  E = "  -0.27%  " },
    @{ Row = 5;  D = "569.32";    DText = $true;  E = "  -0.82%  " },
    @{ Row = 6;  D = "144.06";    DText = $true;  E = "  -3.92%  " },
    @{ Row = 7;  D = $null;       DText = $false; E = "  +0.13%  " },
    @{ Row = 8;  D = $null;       DText = $false; E = "  -1.55%  " },
    @{ Row = 9;  D = "2.438.40";  DText = $false; E = "  -1.95%  " },
    @{ Row = 10; D = $null;       DText = $false; E = "  -4.05%  " },
    @{ Row = 11; D = $null;       DText = $false; E = "  +1.44%  " },
    @{ Row = 12; D = $null;       DText = $false; E = "  -1.73%  " },
    @{ Row = 13; D = $null;       DText = $false; E = "  -2.77%  " },
    @{ Row = 14; D = "26.90";     DText = $true;  E = "  -1.44%  " },
    @{ Row = 15; D = $null;       DText = $false; E = "  -5.39%  " },
    @{ Row = 16; D = "2.884.03";  DText = $false; E = "  -0.88%  " },
    @{ Row = 17; D = "62.427.05"; DText = $false; E = "  -1.49%  " },
    @{ Row = 18; D = "2.425.96";  DText = $false; E = "  -2.56%  " },
    @{ Row = 19; D = "11.19";     DText = $true;  E = "  -3.51%  " },
    @{ Row = 20; D = "7.25";      DText = $true;  E = "  -0.14%  " },
    @{ Row = 21; D = "326.47";    DText = $true;  E = "  -0.74%  " },
    @{ Row = 22; D = "4.16";      DText = $true;  E = "  -2.10%  " },
    @{ Row = 23; D = "2.10";      DText = $true;  E = "  +9.86%  " },
    @{ Row = 24; D = $null;       DText = $false; E = "  +0.19%  " },
    @{ Row = 25; D = "65.24";     DText = $true;  E = "  -3.51%  " },
    @{ Row = 26; D = "621.78";    DText = $true;  E = "  -2.68%  " },
    @{ Row = 27; D = "9.01";      DText = $true;  E = "  +1.77%  " },
    @{ Row = 28; D = "0.0₃0994";  DText = $false; E = "  -5.46%  " },
    @{ Row = 29; D = $null;       DText = $false; E = "  -1.51%  " },
    @{ Row = 30; D = "1.00";      DText = $true;  E = "  +1.22%  " },
    @{ Row = 31; D = "1.48";      DText = $true;  E = "  -2.60%  " },
    @{ Row = 32; D = "8.12";      DText = $true;  E = "  -4.47%  " },
    @{ Row = 33; D = "1.88";      DText = $true;  E = "  -1.94%  " },
    @{ Row = 34; D = $null;       DText = $false; E = "  -4.14%  " },
    @{ Row = 35; D = "5.11";      DText = $true;  E = "  -2.62%  " },
    @{ Row = 36; D = "1.49";      DText = $true;  E = "  -4.59%  " },
    @{ Row = 37; D = "0.999";     DText = $true;  E = "  +0.19%  " },
    @{ Row = 38; D = $null;       DText = $false; E = "  -2.72%  " },
    @{ Row = 39; D = "18.79";     DText = $true;  E = "  -0.92%  " },
    @{ Row = 40; D = "5.33";      DText = $true;  E = "  -3.21%  " },
    @{ Row = 41; D = "147.17";    DText = $true;  E = "  +0.23%  " },
    @{ Row = 42; D = $null;       DText = $false; E = "  -4.86%  " },
    @{ Row = 43; D = "2.56";      DText = $true;  E = "  -2.71%  " },
    @{ Row = 46; D = "146.25";    DText = $true;  E = "  -3.56%  " },
    @{ Row = 47; D = "3.74";      DText = $true;  E = "  -1.21%  " },
    @{ Row = 48; D = "20.66";     DText = $true;  E = "  -2.26%  " },
    @{ Row = 49; D = "0.0529";    DText = $true;  E = "  -4.40%  " },
    @{ Row = 50; D = "0.599";     DText = $true;  E = "  -2.10%  " },
    @{ Row = 51; D = "0.0230";    DText = $true;  E = "  -4.19%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        if ($u.DText) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}

# Rows 44 and 45 swapped coin identities (OKB <-> USDe) with new values.
$ws.Cells.Item(44, 2).Value = "USDe"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$d44 = $ws.Cells.Item(44, 4)
$d44.NumberFormat = "@"
$d44.Value = "0.999"
$ws.Cells.Item(44, 5).Value = "  +0.00%  "

$ws.Cells.Item(45, 2).Value = "OKB"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$d45 = $ws.Cells.Item(45, 4)
$d45.NumberFormat = "@"
$d45.Value = "42.17"
$ws.Cells.Item(45, 5).Value = "  +1.00%  "
